$d = $word.ActiveDocument

# --- 1. Insert the new "post" verbatim line + line break before the
#        first existing "pt_base age" verbatim line in the SourceCode
#        paragraph. ---------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("pt_base age", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Collapse(1)
$insStart = $rng.Start

# 6 leading spaces + the Stata command line (note: the backtick is a
# literal backtick character, escaped for PowerShell with a second
# backtick; the apostrophe is a plain ASCII apostrophe).
$newText = "      . post ``postname' (`"Variable`") (`"Cat level`") (`"Group1`") (`"Group2`") (`"Overall`")"

$rng.InsertBefore($newText)

$insRange = $d.Range($insStart, $insStart + $newText.Length)
$insRange.CharacterStyle = "VerbatimChar"

$brPos = $insStart + $newText.Length
$brRange = $d.Range($brPos, $brPos)
$brRange.InsertBreak(6)

# --- 2. Bump the nsid of abstractNum 990 in numbering.xml. -----------
# (handled separately below, if the runtime exposes list definitions)
